$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

Set-TextValue "D2" "276.18"
Set-TextValue "E2" "-1.22%"
Set-TextValue "D3" "26.73"
Set-TextValue "E3" "-2.24%"
Set-TextValue "D4" "4.876"
Set-TextValue "E4" "1.18%"
Set-TextValue "D5" "0.06327"
Set-TextValue "E5" "-0.15%"
Set-TextValue "D6" "6.910"
Set-TextValue "E6" "-0.46%"
Set-TextValue "D7" "1.321"
Set-TextValue "E7" "39.37%"
Set-TextValue "D8" "0.8756"
Set-TextValue "E8" "-0.44%"
Set-TextValue "D9" "0.1550"
Set-TextValue "E9" "5.41%"
Set-TextValue "D10" "0.05011"
Set-TextValue "E10" "-2.49%"
Set-TextValue "D11" "0.07478"
Set-TextValue "E11" "2.76%"
Set-TextValue "D12" "0.02926"
Set-TextValue "E12" "-6.61%"
Set-TextValue "D13" "0.09051"
Set-TextValue "E13" "-0.18%"
Set-TextValue "D14" "0.001572"
Set-TextValue "E14" "1.27%"
Set-TextValue "D15" "0.0006324"
Set-TextValue "E15" "0.67%"
Set-TextValue "D16" "0.005857"
Set-TextValue "E16" "-2.49%"
Set-TextValue "E17" "0.24%"
Set-TextValue "D18" "3.318"
Set-TextValue "E18" "-1.99%"
Set-TextValue "D19" "2.284"
Set-TextValue "E19" "-0.33%"
Set-TextValue "E21" "1.72%"
Set-TextValue "D22" "3.899"
Set-TextValue "E22" "0.92%"
Set-TextValue "D23" "0.04350"
Set-TextValue "E23" "0.32%"
Set-TextValue "E24" "-0.78%"
Set-TextValue "D25" "0.004209"
Set-TextValue "E25" "-2.03%"
Set-TextValue "D26" "0.0001200"
Set-TextValue "E26" "0.04%"
Set-TextValue "D27" "0.0001671"
Set-TextValue "E27" "-1.03%"
Set-TextValue "D40" "0.04095"
Set-TextValue "E40" "0.25%"
Set-TextValue "D41" "0.006970"
Set-TextValue "E41" "4.21%"
Set-TextValue "D42" "0.1174"
Set-TextValue "E42" "0.94%"
Set-TextValue "D43" "0.002279"
Set-TextValue "E43" "3.68%"
Set-TextValue "D44" "0.01078"
Set-TextValue "E44" "-17.53%"
Set-TextValue "D45" "0.00005282"
Set-TextValue "E45" "1.09%"
$ws.Range("B46").Value = "CoinbaseStockToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextValue "D46" "0.02099"
Set-TextValue "E46" "-6.63%"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextValue "D47" "1.490"
Set-TextValue "E47" "-37.36%"
